$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write plain label/header text (never numeric-looking).
function Set-Text($addr, $text) {
    $ws.Range($addr).Value = $text
}

# Helper: write a value that looks numeric but must be stored as literal text
# (matches the source workbook, where figures like "-0.021" / "1.049***" are
# shared-string text, not numbers). A leading apostrophe forces Excel to treat
# it as text; resetting the style back to "Normal" afterwards clears the
# quote-prefix formatting so the cell keeps the default (unstyled) look.
function Set-TextNumber($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).Style = "Normal"
}

# Clear only the previous 3x3 table's contents (keep the existing cell
# formatting - the bold/bordered/centered style on the header row and the
# label column survives for the cells that are reused).
$ws.Range("A1:C3").ClearContents()

# Header row
Set-Text "A1" "Source"
Set-Text "B1" "C"
Set-Text "C1" "FFR"
Set-Text "D1" "LF"

# Row labels
Set-Text "A2" "C Lag"
Set-Text "A3" "FFR Lag"
Set-Text "A4" "LF Lag"

# Data values (stored as text, same convention as the original workbook)
Set-TextNumber "B2" "-0.46***"
Set-TextNumber "C2" "3.79"
Set-TextNumber "D2" "-6.09"

Set-TextNumber "B3" "-0.01"
Set-TextNumber "C3" "1.6***"
Set-TextNumber "D3" "0.5***"

Set-TextNumber "B4" "0.04*"
Set-TextNumber "C4" "3.53*"
Set-TextNumber "D4" "0.54*"

# Propagate the bold/bordered/centered header style (already present on A1)
# to the newly-added header cell D1 and the newly-added label cell A4.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
